# db test dada added
# Adds a sample DB record in column D of the "in column format" sheet,
# mirroring the vertical key/value layout already used in column A/B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlPasteFormats = -4122

# --- "plain" (General number format, centered) cells -----------------
# First cell establishes the shared style; its format is then copied to
# every other "plain" cell so they all reuse the very same style index.
$ws.Range("D2").HorizontalAlignment = $xlCenter
$ws.Range("D2").Value = 123123

$ws.Range("D2").Copy()
$plainCells = "D3", "D4", "D9", "D10", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D30", "D31"
foreach ($addr in $plainCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$ws.Application.CutCopyMode = $false

$ws.Range("D3").Value = "active"
$ws.Range("D4").Value = "active"
$ws.Range("D9").Value = "scheme2023"
$ws.Range("D10").Value = "active"
$ws.Range("D13").Value = 100000
$ws.Range("D14").Value = 5000
$ws.Range("D15").Value = 4000
$ws.Range("D18").Value = 20
$ws.Range("D19").Value = 20
$ws.Range("D20").Value = 2
$ws.Range("D21").Value = 165000
$ws.Range("D30").Value = "emcashed"
$ws.Range("D31").Value = "ramesh"

# --- date cells (mm-dd-yy number format, centered) --------------------
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("D11").NumberFormat = "mm-dd-yy"
$ws.Range("D11").Value = 45084

$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial($xlPasteFormats)
$ws.Range("D22").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("D12").Value = 45877
$ws.Range("D22").Value = 45115

# --- column width -------------------------------------------------------
$ws.Columns("D").ColumnWidth = 16.26

# --- selection / view -----------------------------------------------
$ws.Activate()
$ws.Range("D13").Select()
